$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 to 45) holds the "Förändrad" (Changed) date, stored as serial date 45172.
# Update it to serial date 45175 (2023-09-06) for all rows.
$ws.Range("C2:C45").Value = 45175
